$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7.800818559179596
$ws.Range("D2").Value = 4.786548219736237
$ws.Range("E2").Value = 11.90286008835687
$ws.Range("F2").Value = 26.50698528322041
$ws.Range("G2").Value = 33.11511128756324
$ws.Range("H2").Value = 14.949068772853
$ws.Range("K2").Value = 10.38421041228618
$ws.Range("M2").Value = 14.20588426489885
$ws.Range("N2").Value = 18.93266549881341
$ws.Range("B3").Value = 7.729303826301262
$ws.Range("D3").Value = 4.783410153023312
$ws.Range("E3").Value = 11.76327144579289
$ws.Range("F3").Value = 26.37109898924724
$ws.Range("G3").Value = 32.81526778124211
$ws.Range("H3").Value = 14.96476793590236
$ws.Range("K3").Value = 9.999814629666608
$ws.Range("M3").Value = 13.96766326380751
$ws.Range("N3").Value = 19.00062134317207
$ws.Range("B4").Value = 7.686933560431236
$ws.Range("D4").Value = 4.781485597484888
$ws.Range("E4").Value = 11.68078601558786
$ws.Range("F4").Value = 26.29575170107471
$ws.Range("G4").Value = 32.64225792513611
$ws.Range("H4").Value = 14.97772152532224
$ws.Range("K4").Value = 9.758431930294078
$ws.Range("M4").Value = 13.82327223646492
$ws.Range("N4").Value = 19.04424496607256
$ws.Range("B5").Value = 7.670073646268931
$ws.Range("D5").Value = 4.78070216771493
$ws.Range("E5").Value = 11.64802111565714
$ws.Range("F5").Value = 26.26710433890891
$ws.Range("G5").Value = 32.57461916748456
$ws.Range("H5").Value = 14.98383187911979
$ws.Range("K5").Value = 9.658881761742458
$ws.Range("M5").Value = 13.76498425307348
$ws.Range("N5").Value = 19.06250075891744
$ws.Range("B6").Value = 7.667299152310518
$ws.Range("D6").Value = 4.780572138979287
$ws.Range("E6").Value = 11.64263285735404
$ws.Range("F6").Value = 26.26247234141389
$ws.Range("G6").Value = 32.56356282420243
$ws.Range("H6").Value = 14.98489667250869
$ws.Range("K6").Value = 9.642285065298728
$ws.Range("M6").Value = 13.75534128879385
$ws.Range("N6").Value = 19.06556108052142
$ws.Range("B7").Value = 7.686704511980355
$ws.Range("D7").Value = 4.781475028109938
$ws.Range("E7").Value = 11.68034065184627
$ws.Range("F7").Value = 26.29535699441794
$ws.Range("G7").Value = 32.64133403551534
$ws.Range("H7").Value = 14.97780056688594
$ws.Range("K7").Value = 9.757093931213237
$ws.Range("M7").Value = 13.82248380131622
$ws.Range("N7").Value = 19.04448922951788
$ws.Range("B8").Value = 7.775852293889021
$ws.Range("D8").Value = 4.78546591821394
$ws.Range("E8").Value = 11.85408464947415
$ws.Range("F8").Value = 26.45846836662861
$ws.Range("G8").Value = 33.00946898705971
$ws.Range("H8").Value = 14.95379301786385
$ws.Range("K8").Value = 10.25287694013361
$ws.Range("M8").Value = 14.12340428319697
$ws.Range("N8").Value = 18.95570351347995
$ws.Range("B9").Value = 7.962019487494559
$ws.Range("D9").Value = 4.793297867190082
$ws.Range("E9").Value = 12.21859132273062
$ws.Range("F9").Value = 26.84129797599048
$ws.Range("G9").Value = 33.81583637981215
$ws.Range("H9").Value = 14.93307737252294
$ws.Range("K9").Value = 11.17774419243782
$ws.Range("M9").Value = 14.7248731993321
$ws.Range("N9").Value = 18.79659185599653
$ws.Range("B10").Value = 8.104509347719024
$ws.Range("D10").Value = 4.799040328053452
$ws.Range("E10").Value = 12.4984207407147
$ws.Range("F10").Value = 27.15911719370006
$ws.Range("G10").Value = 34.45456580621064
$ws.Range("H10").Value = 14.93400269625462
$ws.Range("K10").Value = 11.99780960302421
$ws.Range("M10").Value = 15.16902207850609
$ws.Range("N10").Value = 18.68873992590369
$ws.Range("B11").Value = 8.170308283362292
$ws.Range("D11").Value = 4.801647060689461
$ws.Range("E11").Value = 12.62778757568017
$ws.Range("F11").Value = 27.31118568160639
$ws.Range("G11").Value = 34.75399252426134
$ws.Range("H11").Value = 14.9379383972866
$ws.Range("K11").Value = 12.35115392045079
$ws.Range("M11").Value = 15.37060599040084
$ws.Range("N11").Value = 18.64161952653449
$ws.Range("B12").Value = 8.19534315571571
$ws.Range("D12").Value = 4.802633100122994
$ws.Range("E12").Value = 12.67702704140589
$ws.Range("F12").Value = 27.36980515094304
$ws.Range("G12").Value = 34.86854733631034
$ws.Range("H12").Value = 14.9399342592193
$ws.Range("K12").Value = 12.48210588156086
$ws.Range("M12").Value = 15.44679750736794
$ws.Range("N12").Value = 18.62405401736496
$ws.Range("B13").Value = 8.189946548764105
$ws.Range("D13").Value = 4.802420790994366
$ws.Range("E13").Value = 12.6664120181405
$ws.Range("F13").Value = 27.35713509747467
$ws.Range("G13").Value = 34.84382565783869
$ws.Range("H13").Value = 14.93948193512005
$ws.Range("K13").Value = 12.45403010162874
$ws.Range("M13").Value = 15.43039590395715
$ws.Range("N13").Value = 18.62782472198659
$ws.Range("B14").Value = 8.172365677009578
$ws.Range("D14").Value = 4.801728205835292
$ws.Range("E14").Value = 12.63183377182381
$ws.Range("F14").Value = 27.31598786201312
$ws.Range("G14").Value = 34.76339421755303
$ws.Range("H14").Value = 14.938092467975
$ws.Range("K14").Value = 12.36198465692761
$ws.Range("M14").Value = 15.37687763013436
$ws.Range("N14").Value = 18.64016883807997
$ws.Range("B15").Value = 8.161611613543281
$ws.Range("D15").Value = 4.801303829464911
$ws.Range("E15").Value = 12.61068490734729
$ws.Range("F15").Value = 27.29091739804871
$ws.Range("G15").Value = 34.71427660645703
$ws.Range("H15").Value = 14.93730720602348
$ws.Range("K15").Value = 12.30523232462742
$ws.Range("M15").Value = 15.3440751489168
$ws.Range("N15").Value = 18.64776612652492
$ws.Range("B16").Value = 8.100227010237464
$ws.Range("D16").Value = 4.798869836674643
$ws.Range("E16").Value = 12.49000398822413
$ws.Range("F16").Value = 27.14932629555107
$ws.Range("G16").Value = 34.43516691296381
$ws.Range("H16").Value = 14.93381622998738
$ws.Range("K16").Value = 11.97431949837559
$ws.Range("M16").Value = 15.1558321071811
$ws.Range("N16").Value = 18.69185832183057
$ws.Range("B17").Value = 8.062804344169468
$ws.Range("D17").Value = 4.797375043321812
$ws.Range("E17").Value = 12.41646729709049
$ws.Range("F17").Value = 27.06435367890258
$ws.Range("G17").Value = 34.26613780704144
$ws.Range("H17").Value = 14.93257501596868
$ws.Range("K17").Value = 11.76625312020616
$ws.Range("M17").Value = 15.04017672594865
$ws.Range("N17").Value = 18.71940392161942
$ws.Range("B18").Value = 8.041373376281207
$ws.Range("D18").Value = 4.79651477233095
$ws.Range("E18").Value = 12.37436867496544
$ws.Range("F18").Value = 27.0161874758305
$ws.Range("G18").Value = 34.16975773555516
$ws.Range("H18").Value = 14.93219194701038
$ws.Range("K18").Value = 11.64472456991259
$ws.Range("M18").Value = 14.97361581678463
$ws.Range("N18").Value = 18.73543028719779
$ws.Range("B19").Value = 8.034133977172008
$ws.Range("D19").Value = 4.796223421075769
$ws.Range("E19").Value = 12.36015019295878
$ws.Range("F19").Value = 27.0000020564781
$ws.Range("G19").Value = 34.13727274905484
$ws.Range("H19").Value = 14.93211906104454
$ws.Range("K19").Value = 11.60325925314873
$ws.Range("M19").Value = 14.95107515230908
$ws.Range("N19").Value = 18.74088798606127
$ws.Range("B20").Value = 8.066778531698073
$ws.Range("D20").Value = 4.797534220933935
$ws.Range("E20").Value = 12.42427527747499
$ws.Range("F20").Value = 27.07332620677631
$ws.Range("G20").Value = 34.28404496785518
$ws.Range("H20").Value = 14.93267290239062
$ws.Range("K20").Value = 11.78859419840293
$ws.Range("M20").Value = 15.0524930458629
$ws.Range("N20").Value = 18.7164527290047
$ws.Range("B21").Value = 8.177526576167061
$ws.Range("D21").Value = 4.801931666206686
$ws.Range("E21").Value = 12.64198379913611
$ws.Range("F21").Value = 27.32804608025609
$ws.Range("G21").Value = 34.78698802891545
$ws.Range("H21").Value = 14.93848687027864
$ws.Range("K21").Value = 12.38909815112842
$ws.Range("M21").Value = 15.39260175337602
$ws.Range("N21").Value = 18.63653553834282
$ws.Range("B22").Value = 8.250585604056578
$ws.Range("D22").Value = 4.804799267551667
$ws.Range("E22").Value = 12.78571075198878
$ws.Range("F22").Value = 27.50052834010027
$ws.Range("G22").Value = 35.12245021459612
$ws.Range("H22").Value = 14.94523283890398
$ws.Range("K22").Value = 12.76493453409492
$ws.Range("M22").Value = 15.61400979450389
$ws.Range("N22").Value = 18.58592462706065
$ws.Range("B23").Value = 8.211537973071444
$ws.Range("D23").Value = 4.803269449166228
$ws.Range("E23").Value = 12.70888442447538
$ws.Range("F23").Value = 27.40793626775184
$ws.Range("G23").Value = 34.94282484109074
$ws.Range("H23").Value = 14.94136287887052
$ws.Range("K23").Value = 12.56586963398338
$ws.Range("M23").Value = 15.49594480137544
$ws.Range("N23").Value = 18.61278885377922
$ws.Range("B24").Value = 8.064981538707842
$ws.Range("D24").Value = 4.797462259468525
$ws.Range("E24").Value = 12.42074472968513
$ws.Range("F24").Value = 27.06926758772712
$ws.Range("G24").Value = 34.27594664960938
$ws.Range("H24").Value = 14.93262761831928
$ws.Range("K24").Value = 11.77849973996735
$ws.Range("M24").Value = 15.04692504687428
$ws.Range("N24").Value = 18.71778637108908
$ws.Range("B25").Value = 7.910564188879787
$ws.Range("D25").Value = 4.791180082755543
$ws.Range("E25").Value = 12.11768688274622
$ws.Range("F25").Value = 26.73116814475443
$ws.Range("G25").Value = 33.5891680836242
$ws.Range("H25").Value = 14.93585025092749
$ws.Range("K25").Value = 10.93246211313531
$ws.Range("M25").Value = 14.56143170211663
$ws.Range("N25").Value = 18.83804012943622
